# "adding averages and more checks"
#
# - Training Dashboard: "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I)
#   refreshed for rows 3-15 (the LAST UPDATE date moved from 08-Sep-2025 to
#   16-Sep-2025, eight days later, so PERIOD TO EXPIRE drops by 8 for every row).
# - Exam Dashboard: widen the COMMENTS column and make the per-row comment
#   more descriptive ("OK" -> "date is valid").
# - Header rows get white bold text so they read clearly against the dark
#   blue fill (previously the bold header font had no explicit color).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Training Dashboard: refresh PERIOD TO EXPIRE / LAST UPDATE
# ---------------------------------------------------------------------------
$training = $wb.Worksheets.Item("Training Dashboard")

$updates = @(
    @{ Row = 3;  Period = 426;  },
    @{ Row = 4;  Period = 267;  },
    @{ Row = 5;  Period = 358;  },
    @{ Row = 6;  Period = 253;  },
    @{ Row = 7;  Period = 219;  },
    @{ Row = 8;  Period = 266;  },
    @{ Row = 9;  Period = 400;  },
    @{ Row = 10; Period = 335;  },
    @{ Row = 11; Period = 350;  },
    @{ Row = 12; Period = 357;  },
    @{ Row = 13; Period = -22;  },
    @{ Row = 14; Period = -103; },
    @{ Row = 15; Period = 155;  }
)

foreach ($u in $updates) {
    $r = $u.Row

    # Column H ("PERIOD TO EXPIRE") is a plain number, safe to assign directly.
    $training.Cells.Item($r, 8).Value = $u.Period

    # Column I ("LAST UPDATE") holds the date as literal text (General
    # format). A bare assignment of a date-shaped string gets auto-parsed
    # into a real date serial by Excel, so prefix with a leading apostrophe
    # -- exactly like typing it in by hand -- to force literal text entry
    # and keep it the string "16-Sep-2025" (General format, unchanged).
    $training.Cells.Item($r, 9).Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# Exam Dashboard: wider COMMENTS column + clearer per-row remark
# ---------------------------------------------------------------------------
$exam = $wb.Worksheets.Item("Exam Dashboard")

$exam.Columns.Item(5).ColumnWidth = 14.166666666666666

foreach ($r in 3..8) {
    $exam.Cells.Item($r, 5).Value = "date is valid"
}

# ---------------------------------------------------------------------------
# Header styling: bold white text on the dark-blue header band, and the big
# title keeps bold but shrinks to the normal body size.
# ---------------------------------------------------------------------------
$white = 16777215

$training.Range("A2:K2").Font.Color = $white
$exam.Range("A2:G2").Font.Color = $white

foreach ($ws in @($training, $exam)) {
    $title = $ws.Range("A1")
    $title.Font.Color = $white
    $title.Font.Size = 11
    $title.Font.Bold = $true
}

Write-Host "done"
